# Loan RBI, Variable Instalments
# Inserts a new (blank) column into the "Repayment schedule" sheet between
# the existing "In Advance" (M) and "Late" (old N) columns, shifting the
# trailing columns (Late / Paid Date+Disbursement / Outstanding) one to the
# right, then makes "Repayment schedule" the active sheet/selection.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a whole new column at N; existing N (Late), O (Paid Date/heading)
# and P (Outstanding) shift right to O, P, Q respectively.
$ws.Columns("N").Insert() | Out-Null

# Give the freshly inserted column the same width as its left-hand
# neighbour (column M, "In Advance") instead of the narrow default.
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Make "Repayment schedule" the active sheet and select cell O6 on it.
$ws.Activate() | Out-Null
$ws.Range("O6").Select() | Out-Null
